$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A38").Value = 54
$ws.Range("A45").Value = 55

$ws.Range("A77").Value = 56
$ws.Range("C77").Value = "ევროპული საქართველო"

$ws.Range("A78").Value = 57
$ws.Range("C78").Value = "თავისუფლება ზვიად გამსახურდიას გზა"

$ws.Range("A79").Value = 58
$ws.Range("C79").Value = "ნეიტრალური საქართველო"

$ws.Range("A80").Value = 59
$ws.Range("C80").Value = "სრულიად საქართველოს რადიკალ-დემოკრატთა ნაციონალური პარტია"

$ws.Range("A81").Value = 60
$ws.Range("C81").Value = "მოქალაქეთა პოლიტიკური გაერთიანება სახალხო ხელისუფლება"
